$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")

# Row 41
$ws.Range("H41").Value = 29415574
$ws.Range("J41").Value = 45459716
$ws.Range("L41").Value = 45459716
$ws.Range("N41").Value = -45460596
# Row 51
$ws.Range("H51").Value = 7578.25
$ws.Range("J51").Value = 6771.3335
$ws.Range("L51").Value = 6771.3335
$ws.Range("N51").Value = -7739.3335
# Row 64
$ws.Range("H64").Value = 8087.875
$ws.Range("J64").Value = 9000.5
$ws.Range("L64").Value = 9000.5
$ws.Range("N64").Value = -9496.5
# Row 67
$ws.Range("H67").Value = 8087.875
$ws.Range("J67").Value = 9000.5
$ws.Range("L67").Value = 9000.5
$ws.Range("N67").Value = -10716.5
# Row 70
$ws.Range("H70").Value = 56673.055
$ws.Range("J70").Value = 92216.27
$ws.Range("L70").Value = 276648.81
$ws.Range("N70").Value = -277188.81
# Row 73
$ws.Range("H73").Value = 56673.055
$ws.Range("J73").Value = 92216.27
$ws.Range("L73").Value = 276648.81
$ws.Range("N73").Value = -278520.81
# Row 80
$ws.Range("H80").Value = 1206.5
$ws.Range("I80").Value = 989.5
$ws.Range("J80").Value = 1351.1666
$ws.Range("K80").Value = 2968.5
$ws.Range("L80").Value = 4053.4998
$ws.Range("M80").Value = -1970.5
$ws.Range("N80").Value = -6049.4998
# Row 83
$ws.Range("H83").Value = 1206.5
$ws.Range("I83").Value = 989.5
$ws.Range("J83").Value = 1351.1666
$ws.Range("K83").Value = 8905.5
$ws.Range("L83").Value = 12160.4994
$ws.Range("M83").Value = -3913.5
$ws.Range("N83").Value = -22144.4994
# Row 100
$ws.Range("H100").Value = 4885.357
$ws.Range("J100").Value = 5399.5454
$ws.Range("L100").Value = 5399.5454
$ws.Range("N100").Value = -6481.5454
# Row 106
$ws.Range("H106").Value = 1742.7273
$ws.Range("I106").Value = 1663.3334
$ws.Range("J106").Value = 1838
$ws.Range("K106").Value = 1663.3334
$ws.Range("L106").Value = 1838
$ws.Range("M106").Value = -1032.3334
$ws.Range("N106").Value = -3100
# Row 113
$ws.Range("H113").Value = 15279.154
$ws.Range("I113").Value = 13785.429
$ws.Range("J113").Value = 17021.834
$ws.Range("K113").Value = 13785.429
$ws.Range("L113").Value = 17021.834
$ws.Range("M113").Value = -10531.429
$ws.Range("N113").Value = -23529.834
# Row 132
$ws.Range("H132").Value = 5093
$ws.Range("I132").Value = 5193.16
$ws.Range("K132").Value = 15579.48
$ws.Range("M132").Value = -13049.48
# Row 137
$ws.Range("H137").Value = 3464.56
$ws.Range("I137").Value = 2112
$ws.Range("K137").Value = 6336
$ws.Range("M137").Value = -3786

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")

# Row 2
$ws.Range("H2").Value = 5832.8887
$ws.Range("I2").Value = 5337
$ws.Range("K2").Value = 5337
$ws.Range("M2").Value = -5224
# Row 32
$ws.Range("H32").Value = 2088.7097
$ws.Range("I32").Value = 2116.4753
$ws.Range("K32").Value = 2116.4753
$ws.Range("M32").Value = -1829.4753
# Row 94
$ws.Range("H94").Value = 22500
$ws.Range("J94").Value = 22500
$ws.Range("L94").Value = 22500
$ws.Range("N94").Value = -24302
# Row 102
$ws.Range("H102").Value = 2096.7856
$ws.Range("I102").Value = 1877.6
$ws.Range("K102").Value = 1877.6
$ws.Range("M102").Value = -255.5999999999999
# Row 109
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("N109").Value = 0
# Row 116
$ws.Range("H116").Value = 5832.8887
$ws.Range("I116").Value = 5337
$ws.Range("K116").Value = 5337
$ws.Range("M116").Value = -3043
# Row 132
$ws.Range("H132").Value = 9412.869000000001
$ws.Range("I132").Value = 3422.8462
$ws.Range("K132").Value = 10268.5386
$ws.Range("M132").Value = -7738.5386

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")

# Row 3
$ws.Range("H3").Value = 5832.8887
$ws.Range("I3").Value = 5337
$ws.Range("K3").Value = 5337
$ws.Range("M3").Value = -5223
# Row 26
$ws.Range("H26").Value = 21123.666
$ws.Range("I26").Value = 16685.5
$ws.Range("J26").Value = 30000
$ws.Range("K26").Value = 16685.5
$ws.Range("L26").Value = 30000
$ws.Range("M26").Value = -16393.5
$ws.Range("N26").Value = -30584
# Row 106
$ws.Range("H106").Value = 9999
$ws.Range("J106").Value = 9999
$ws.Range("L106").Value = 9999
$ws.Range("N106").Value = -12523

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")

# Row 22
$ws.Range("H22").Value = 725
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
# Row 132
$ws.Range("H132").Value = 4702.143
$ws.Range("I132").Value = 4193.5
$ws.Range("K132").Value = 12580.5
$ws.Range("M132").Value = -10050.5

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")

# Row 12
$ws.Range("H12").Value = 207.6
$ws.Range("J12").Value = 232.3077
$ws.Range("L12").Value = 696.9231
$ws.Range("N12").Value = -1042.9231
# Row 33
$ws.Range("H33").Value = 2645520.8
$ws.Range("I33").Value = 3086439.2
$ws.Range("J33").Value = 9
$ws.Range("K33").Value = 18518635.2
$ws.Range("L33").Value = 54
$ws.Range("M33").Value = -18518352.2
$ws.Range("N33").Value = -620
# Row 128
$ws.Range("H128").Value = 144496
$ws.Range("I128").Value = 144496
$ws.Range("K128").Value = 433488
$ws.Range("M128").Value = -428508

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")

# Row 13
$ws.Range("H13").Value = 1572.8889
$ws.Range("I13").Value = 1675
$ws.Range("K13").Value = 1675
$ws.Range("M13").Value = -1536
# Row 41
$ws.Range("H41").Value = 1332.25
$ws.Range("I41").Value = 976.3333
$ws.Range("K41").Value = 976.3333
$ws.Range("M41").Value = -621.3333
# Row 42
$ws.Range("H42").Value = 49999.332
$ws.Range("J42").Value = 49999.332
$ws.Range("L42").Value = 49999.332
$ws.Range("N42").Value = -50969.332
# Row 115
$ws.Range("H115").Value = 49999.332
$ws.Range("J115").Value = 49999.332
$ws.Range("L115").Value = 49999.332
$ws.Range("N115").Value = -52349.332
# Row 133
$ws.Range("H133").Value = 69983.336
$ws.Range("J133").Value = 69983.336
$ws.Range("L133").Value = 69983.336
$ws.Range("N133").Value = -80103.336

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")

# Row 7
$ws.Range("H7").Value = 423727.72
$ws.Range("I7").Value = 776774.75
$ws.Range("J7").Value = 6490.273
$ws.Range("K7").Value = 776774.75
$ws.Range("L7").Value = 6490.273
$ws.Range("M7").Value = -776662.75
$ws.Range("N7").Value = -6714.273
# Row 75
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("M75").ClearContents()
# Row 78
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("M78").ClearContents()
# Row 126
$ws.Range("H126").Value = 423727.72
$ws.Range("I126").Value = 776774.75
$ws.Range("J126").Value = 6490.273
$ws.Range("K126").Value = 2330324.25
$ws.Range("L126").Value = 19470.819
$ws.Range("M126").Value = -2327854.25
$ws.Range("N126").Value = -24410.819
# Row 132
$ws.Range("H132").Value = 5742.857
$ws.Range("I132").Value = 3771.4285
$ws.Range("J132").Value = 7714.2856
$ws.Range("K132").Value = 11314.2855
$ws.Range("L132").Value = 23142.8568
$ws.Range("M132").Value = -8784.2855
$ws.Range("N132").Value = -28202.8568
# Row 136
$ws.Range("H136").Value = 3993.5
$ws.Range("I136").Value = 2968.7
$ws.Range("K136").Value = 8906.099999999999
$ws.Range("M136").Value = -6356.099999999999

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")

# Row 47
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("N47").Value = 0
# Row 121
$ws.Range("H121").Value = 40737.75
$ws.Range("J121").Value = 40737.75
$ws.Range("L121").Value = 40737.75
$ws.Range("N121").Value = -44231.75
# Row 126
$ws.Range("H126").Value = 4127.9375
$ws.Range("I126").Value = 3262.375
$ws.Range("J126").Value = 4993.5
$ws.Range("K126").Value = 9787.125
$ws.Range("L126").Value = 14980.5
$ws.Range("M126").Value = -7317.125
$ws.Range("N126").Value = -19920.5
# Row 132
$ws.Range("H132").Value = 37966.367
$ws.Range("I132").Value = 2792.923
$ws.Range("K132").Value = 8378.769
$ws.Range("M132").Value = -5848.769
# Row 136
$ws.Range("H136").Value = 318941.72
$ws.Range("I136").Value = 335694.7
$ws.Range("K136").Value = 1007084.1
$ws.Range("M136").Value = -1004534.1
